# Auto-generated edit script applying numeric corrections to the
# Pandaemonium_Profits "Sheets" tabs (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 400.07693
$ws.Range("I55").Value = 193.625
$ws.Range("J55").Value = 730.4
$ws.Range("K55").Value = 193.625
$ws.Range("L55").Value = 730.4
$ws.Range("M55").Value = 20.375
$ws.Range("N55").Value = -1158.4

$ws.Range("H130").Value = 79800
$ws.Range("J130").Value = 79800
$ws.Range("L130").Value = 79800
$ws.Range("N130").Value = -89840

$ws.Range("H137").Value = 5631.107
$ws.Range("I137").Value = 5928.4287
$ws.Range("J137").Value = 5333.7856
$ws.Range("K137").Value = 17785.2861
$ws.Range("L137").Value = 16001.3568
$ws.Range("M137").Value = -15235.2861
$ws.Range("N137").Value = -21101.3568

$ws.Range("H138").Value = 4182.1113
$ws.Range("J138").Value = 4864.5483
$ws.Range("L138").Value = 14593.6449
$ws.Range("N138").Value = -24873.6449

$ws.Range("H139").Value = 57563.5
$ws.Range("J139").Value = 57563.5
$ws.Range("L139").Value = 57563.5
$ws.Range("N139").Value = -67843.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9183.235000000001
$ws.Range("I61").Value = 9438.538
$ws.Range("J61").Value = 8353.5
$ws.Range("K61").Value = 9438.538
$ws.Range("L61").Value = 8353.5
$ws.Range("M61").Value = -9226.538
$ws.Range("N61").Value = -8777.5

$ws.Range("H74").Value = 6269.107
$ws.Range("I74").Value = 2613.8
$ws.Range("J74").Value = 10486.77
$ws.Range("K74").Value = 2613.8
$ws.Range("L74").Value = 10486.77
$ws.Range("M74").Value = -1739.8
$ws.Range("N74").Value = -12234.77

$ws.Range("H77").Value = 6269.107
$ws.Range("I77").Value = 2613.8
$ws.Range("J77").Value = 10486.77
$ws.Range("K77").Value = 13069
$ws.Range("L77").Value = 52433.85000000001
$ws.Range("M77").Value = -8701
$ws.Range("N77").Value = -61169.85000000001

$ws.Range("H132").Value = 3715.4707
$ws.Range("I132").Value = 3313.6667
$ws.Range("J132").Value = 4679.8
$ws.Range("K132").Value = 9941.000100000001
$ws.Range("L132").Value = 14039.4
$ws.Range("M132").Value = -7411.000100000001
$ws.Range("N132").Value = -19099.4

$ws.Range("H136").Value = 9183.235000000001
$ws.Range("I136").Value = 9438.538
$ws.Range("J136").Value = 8353.5
$ws.Range("K136").Value = 28315.614
$ws.Range("L136").Value = 25060.5
$ws.Range("M136").Value = -25765.614
$ws.Range("N136").Value = -30160.5

$ws.Range("H137").Value = 56347.25
$ws.Range("J137").Value = 56347.25
$ws.Range("L137").Value = 56347.25
$ws.Range("N137").Value = -66547.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2525.9167
$ws.Range("I107").Value = 2482.818
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 2482.818
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -562.8180000000002
$ws.Range("N107").Value = -6840

$ws.Range("H134").Value = 8523.076999999999
$ws.Range("I134").Value = 7865.1113
$ws.Range("J134").Value = 10003.5
$ws.Range("K134").Value = 23595.3339
$ws.Range("L134").Value = 30010.5
$ws.Range("M134").Value = -21060.3339
$ws.Range("N134").Value = -35080.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 66
$ws.Range("I7").Value = 62.5
$ws.Range("K7").Value = 62.5
$ws.Range("M7").Value = 50.5

$ws.Range("H31").Value = 3917.4565
$ws.Range("I31").Value = 5496.091
$ws.Range("J31").Value = 2470.375
$ws.Range("K31").Value = 5496.091
$ws.Range("L31").Value = 2470.375
$ws.Range("M31").Value = -5201.091
$ws.Range("N31").Value = -3060.375

$ws.Range("H34").Value = 3917.4565
$ws.Range("I34").Value = 5496.091
$ws.Range("J34").Value = 2470.375
$ws.Range("K34").Value = 5496.091
$ws.Range("L34").Value = 2470.375
$ws.Range("M34").Value = -5294.091
$ws.Range("N34").Value = -2874.375

$ws.Range("H62").Value = 3854.1428
$ws.Range("I62").Value = 3854.1428
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3854.1428
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3230.1428
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 3854.1428
$ws.Range("I65").Value = 3854.1428
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 19270.714
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -16150.714
$ws.Range("N65").ClearContents()

$ws.Range("H99").Value = 3342.7144
$ws.Range("I99").Value = 1399.8334
$ws.Range("K99").Value = 1399.8334
$ws.Range("M99").Value = 98.16660000000002

$ws.Range("H107").Value = 387.2143
$ws.Range("I107").Value = 385.16666
$ws.Range("J107").Value = 399.5
$ws.Range("K107").Value = 385.16666
$ws.Range("L107").Value = 399.5
$ws.Range("M107").Value = 1534.83334
$ws.Range("N107").Value = -4239.5

$ws.Range("H126").Value = 3342.7144
$ws.Range("I126").Value = 1399.8334
$ws.Range("K126").Value = 4199.5002
$ws.Range("M126").Value = -1729.5002

$ws.Range("H132").Value = 3239.0857
$ws.Range("I132").Value = 2703.7932
$ws.Range("J132").Value = 5826.3335
$ws.Range("K132").Value = 8111.3796
$ws.Range("L132").Value = 17479.0005
$ws.Range("M132").Value = -5581.3796
$ws.Range("N132").Value = -22539.0005

$ws.Range("H134").Value = 2048.347
$ws.Range("I134").Value = 1688.1892
$ws.Range("J134").Value = 3158.8333
$ws.Range("K134").Value = 5064.5676
$ws.Range("L134").Value = 9476.499899999999
$ws.Range("M134").Value = -2529.5676
$ws.Range("N134").Value = -14546.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2676.1
$ws.Range("I109").Value = 1660.25
$ws.Range("J109").Value = 3353.3333
$ws.Range("K109").Value = 4980.75
$ws.Range("L109").Value = 10059.9999
$ws.Range("M109").Value = -3940.75
$ws.Range("N109").Value = -12139.9999

$ws.Range("H113").Value = 594.1774
$ws.Range("I113").Value = 575.26086
$ws.Range("J113").Value = 648.5625
$ws.Range("K113").Value = 1725.78258
$ws.Range("L113").Value = 1945.6875
$ws.Range("M113").Value = 444.2174199999999
$ws.Range("N113").Value = -6285.6875

$ws.Range("H131").Value = 1078.3088
$ws.Range("J131").Value = 1131.6936
$ws.Range("L131").Value = 3395.0808
$ws.Range("N131").Value = -13475.0808

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 51430
$ws.Range("J32").Value = 51430
$ws.Range("L32").Value = 51430
$ws.Range("N32").Value = -52022

$ws.Range("H42").Value = 59055.715
$ws.Range("J42").Value = 59055.715
$ws.Range("L42").Value = 59055.715
$ws.Range("N42").Value = -60025.715

$ws.Range("H107").Value = 399.25
$ws.Range("I107").Value = 119.8
$ws.Range("J107").Value = 865
$ws.Range("K107").Value = 119.8
$ws.Range("L107").Value = 865
$ws.Range("M107").Value = 1800.2
$ws.Range("N107").Value = -4705

$ws.Range("H113").Value = 2183.2778
$ws.Range("I113").Value = 2075
$ws.Range("J113").Value = 2399.8333
$ws.Range("K113").Value = 2075
$ws.Range("L113").Value = 2399.8333
$ws.Range("M113").Value = 95
$ws.Range("N113").Value = -6739.8333

$ws.Range("H115").Value = 59055.715
$ws.Range("J115").Value = 59055.715
$ws.Range("L115").Value = 59055.715
$ws.Range("N115").Value = -61405.715

$ws.Range("H126").Value = 2858.8235
$ws.Range("I126").Value = 1950
$ws.Range("J126").Value = 3666.6667
$ws.Range("K126").Value = 5850
$ws.Range("L126").Value = 11000.0001
$ws.Range("M126").Value = -3380
$ws.Range("N126").Value = -15940.0001

$ws.Range("H132").Value = 2906.3635
$ws.Range("I132").Value = 1871.5
$ws.Range("J132").Value = 5666
$ws.Range("K132").Value = 5614.5
$ws.Range("L132").Value = 16998
$ws.Range("M132").Value = -3084.5
$ws.Range("N132").Value = -22058

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 6923.923
$ws.Range("I93").Value = 6909.364
$ws.Range("K93").Value = 6909.364
$ws.Range("M93").Value = -5661.364

$ws.Range("H125").Value = 68000
$ws.Range("J125").Value = 68000
$ws.Range("L125").Value = 68000
$ws.Range("N125").Value = -77840

$ws.Range("H132").Value = 3008.1667
$ws.Range("I132").Value = 2113.7
$ws.Range("J132").Value = 3647.0715
$ws.Range("K132").Value = 6341.099999999999
$ws.Range("L132").Value = 10941.2145
$ws.Range("M132").Value = -3811.099999999999
$ws.Range("N132").Value = -16001.2145

$ws.Range("H136").Value = 6515.5713
$ws.Range("I136").Value = 6518.647
$ws.Range("J136").Value = 6502.5
$ws.Range("K136").Value = 19555.941
$ws.Range("L136").Value = 19507.5
$ws.Range("M136").Value = -17005.941
$ws.Range("N136").Value = -24607.5

$ws.Range("H137").Value = 50162.57
$ws.Range("I137").Value = 29390
$ws.Range("J137").Value = 58471.6
$ws.Range("K137").Value = 29390
$ws.Range("L137").Value = 58471.6
$ws.Range("M137").Value = -24290
$ws.Range("N137").Value = -68671.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2824.6
$ws.Range("I107").Value = 1110.8889
$ws.Range("J107").Value = 4226.727
$ws.Range("K107").Value = 3332.6667
$ws.Range("L107").Value = 12680.181
$ws.Range("M107").Value = -1412.6667
$ws.Range("N107").Value = -16520.181

$ws.Range("H126").Value = 1240
$ws.Range("I126").Value = 1283.2727
$ws.Range("J126").Value = 1160.6666
$ws.Range("K126").Value = 3849.8181
$ws.Range("L126").Value = 3481.9998
$ws.Range("M126").Value = -1379.8181
$ws.Range("N126").Value = -8421.9998

$ws.Range("H132").Value = 3421.238
$ws.Range("I132").Value = 3680.2307
$ws.Range("J132").Value = 3000.375
$ws.Range("K132").Value = 11040.6921
$ws.Range("L132").Value = 9001.125
$ws.Range("M132").Value = -8510.6921
$ws.Range("N132").Value = -14061.125

$ws.Range("H136").Value = 6396.0605
$ws.Range("I136").Value = 2619.6155
$ws.Range("J136").Value = 8850.75
$ws.Range("K136").Value = 7858.8465
$ws.Range("L136").Value = 26552.25
$ws.Range("M136").Value = -5308.8465
$ws.Range("N136").Value = -31652.25
